{"js": "// Replace the date line and every \"NNN\u00f7N=\" exercise cell, in document\n// order, with the values from the target revision. The mapping below is\n// positional (old text -> new text) and walks the document's paragraphs\n// (title paragraph + every table-cell paragraph, blank cells included) in\n// natural reading order, so it is robust to the many structurally-identical\n// blank paragraphs that separate the populated rows.\nconst replacements = [\n  [\"2025-08-23 Saturday\", \"2025-08-24 Sunday\"],\n  [\"508\u00f78=\", \"342\u00f72=\"],\n  [\"469\u00f75=\", \"859\u00f73=\"],\n  [\"446\u00f78=\", \"762\u00f73=\"],\n  [\"259\u00f77=\", \"682\u00f78=\"],\n  [\"505\u00f74=\", \"705\u00f75=\"],\n  [\"652\u00f77=\", \"631\u00f77=\"],\n  [\"877\u00f73=\", \"585\u00f78=\"],\n  [\"926\u00f76=\", \"428\u00f79=\"],\n  [\"785\u00f77=\", \"728\u00f74=\"],\n  [\"644\u00f76=\", \"869\u00f72=\"],\n  [\"476\u00f73=\", \"182\u00f72=\"],\n  [\"766\u00f75=\", \"156\u00f77=\"],\n  [\"288\u00f74=\", \"340\u00f79=\"],\n  [\"843\u00f74=\", \"185\u00f76=\"],\n  [\"498\u00f72=\", \"539\u00f75=\"],\n  [\"220\u00f74=\", \"574\u00f79=\"],\n  [\"323\u00f72=\", \"145\u00f76=\"],\n  [\"821\u00f78=\", \"622\u00f72=\"],\n  [\"182\u00f74=\", \"472\u00f74=\"],\n  [\"698\u00f74=\", \"458\u00f72=\"],\n  [\"505\u00f77=\", \"771\u00f76=\"],\n  [\"948\u00f74=\", \"711\u00f76=\"],\n  [\"925\u00f79=\", \"941\u00f74=\"],\n  [\"632\u00f74=\", \"973\u00f72=\"],\n  [\"918\u00f73=\", \"607\u00f73=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet cursor = 0;\nfor (let i = 0; i < paragraphs.items.length && cursor < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const [oldText, newText] = replacements[cursor];\n  if (para.text === oldText) {\n    para.getRange().insertText(newText, \"Replace\");\n    cursor++;\n  }\n}\n\nawait context.sync();\n\nif (cursor !== replacements.length) {\n  throw new Error(\n    `Only replaced ${cursor} of ${replacements.length} expected text runs.`\n  );\n}\n", "ps1": "# Replace the date line and every \"NNN\u00f7N=\" exercise cell with the values\n# from the target revision. Each old value occurs exactly once in the\n# document (title paragraph + 25 table-cell paragraphs), so a simple\n# find-and-replace-one pass per pair, run against the whole document\n# content, reproduces the diff without depending on row/column indices.\n#\n# NOTE: deliberately uses string interpolation / literal string args (never\n# \"digits\" + $char -> PowerShell treats that as numeric addition once the\n# left-hand side parses as a number, silently corrupting values such as\n# \"342\" + [char]247).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-23 Saturday\", \"2025-08-24 Sunday\"),\n    @(\"508\u00f78=\", \"342\u00f72=\"),\n    @(\"469\u00f75=\", \"859\u00f73=\"),\n    @(\"446\u00f78=\", \"762\u00f73=\"),\n    @(\"259\u00f77=\", \"682\u00f78=\"),\n    @(\"505\u00f74=\", \"705\u00f75=\"),\n    @(\"652\u00f77=\", \"631\u00f77=\"),\n    @(\"877\u00f73=\", \"585\u00f78=\"),\n    @(\"926\u00f76=\", \"428\u00f79=\"),\n    @(\"785\u00f77=\", \"728\u00f74=\"),\n    @(\"644\u00f76=\", \"869\u00f72=\"),\n    @(\"476\u00f73=\", \"182\u00f72=\"),\n    @(\"766\u00f75=\", \"156\u00f77=\"),\n    @(\"288\u00f74=\", \"340\u00f79=\"),\n    @(\"843\u00f74=\", \"185\u00f76=\"),\n    @(\"498\u00f72=\", \"539\u00f75=\"),\n    @(\"220\u00f74=\", \"574\u00f79=\"),\n    @(\"323\u00f72=\", \"145\u00f76=\"),\n    @(\"821\u00f78=\", \"622\u00f72=\"),\n    @(\"182\u00f74=\", \"472\u00f74=\"),\n    @(\"698\u00f74=\", \"458\u00f72=\"),\n    @(\"505\u00f77=\", \"771\u00f76=\"),\n    @(\"948\u00f74=\", \"711\u00f76=\"),\n    @(\"925\u00f79=\", \"941\u00f74=\"),\n    @(\"632\u00f74=\", \"973\u00f72=\"),\n    @(\"918\u00f73=\", \"607\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
